$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "on" flag for this sequence step from Hora 0 (column E) to Hora 18 (column W)
# Clear column E (Hora 0) for rows 2-13
$ws.Range("E2:E13").Value = 0

# Set column W (Hora 18) rows 2-4 to 1 (rows 5-13 were already 1)
$ws.Range("W2:W4").Value = 1

# Update the active selection to match the saved workbook state
$ws.Range("G11").Select()
